$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.66%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.48%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.087"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.88%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05597"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.10%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.479"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.19%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8144"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.06%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8443"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.12%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06984"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.45%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.49%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09381"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.25%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.001513"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.17%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.006251"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.51%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.608"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.13%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.020"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.02%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.74%"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.009992"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1,573.57%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3112"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.32%"
$ws.Range("B19").Value = "WazirX"
$ws.Range("C19").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.1333"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.38%"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03198"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.65%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1275"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.39%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.742"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.02%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04655"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.69%"
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1376"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.41%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001249"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.60%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004571"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "6.55%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009604"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-1.00%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001939"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-0.08%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.18%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006153"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.61%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1054"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.16%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002618"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.09%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008062"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.28%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005398"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.91%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.03%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1450"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-19.43%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002426"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "20.35%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
